$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Range("D19").Value = "TIMESTAMP"
$ws.Range("D21").Value = "TIMESTAMP"
$ws.Range("D19").Select()
